# Updated export from Jun's tool
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet: update / remove notes sentences about CCS and process
# emissions policy cost allocation.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A21").Value = "to cover the allocation of capital and OM expensies for industry energy efficiency and CCS."
$wsAbout.Range("A22").ClearContents()
$wsAbout.Range("A23").Value = "A different input variable governs the breakdown of revenues due to Process Emissions policies."
$wsAbout.Range("A24").ClearContents()

# ---------------------------------------------------------------------------
# "SoCaOMSbRIC" sheet: split the combined "ISIC 20T21" column into two
# separate columns, "ISIC 20" and "ISIC 21".
# ---------------------------------------------------------------------------
$wsShare = $wb.Worksheets.Item("SoCaOMSbRIC")

$wsShare.Columns("K").Insert()
$wsShare.Range("K1").Value = "ISIC 20"
$wsShare.Range("L1").Value = "ISIC 21"
$wsShare.Range("K2").Value = 0
